$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list update (prices/volumes refreshed; Frax/EnergySwap rows swapped)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.447.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.107.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.75%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9977"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5252"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4351"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08862"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.87%  "

$ws.Range("E11").Value = "  +2.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.091.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.739"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.47%  "

$ws.Range("E15").Value = "  +4.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9985"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001129"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06633"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.351"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.495.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("E24").Value = "  +5.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.325"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.337.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.596"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.208"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.712"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +25.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1073"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.193"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.919"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02583"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06720"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.488"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2268"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6828"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.252"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.87%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.66%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9965"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6384"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.220"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.613"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.253"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.200"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.74%  "
